$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the "Ajo" (garlic) sheet. In the
# source data the rows are ordered most-recent-first, so the new record is
# inserted right above the existing row 39, pushing it and every row below
# it (old rows 39-152) down by one (new rows 40-153).
$ws.Rows.Item(39).Insert()

# Populate the freshly inserted row 39 with the new record's data.
$ws.Cells.Item(39, 1).Value = 7
$ws.Cells.Item(39, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value = "Ñuble"
$ws.Cells.Item(39, 4).Value = 44525
$ws.Cells.Item(39, 5).Value = 16
$ws.Cells.Item(39, 6).Value = 100112003
$ws.Cells.Item(39, 7).Value = "Ajo"
$ws.Cells.Item(39, 8).Value = "Chino"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 100
$ws.Cells.Item(39, 11).Value = 17000
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = 17500
$ws.Cells.Item(39, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(39, 15).Value = "China"
$ws.Cells.Item(39, 16).Value = 1750
$ws.Cells.Item(39, 17).Value = 10
$ws.Cells.Item(39, 18).Value = "Hortaliza"
